# Update the "ID Competição" column (B) for all data rows (2-73)
# from 48 to 248, per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B73").Value = 248
